$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-10 (row 23)
$ws.Range("B23").Value = 6329
$ws.Range("C23").Value = 1001
$ws.Range("D23").Value = 5923112
$ws.Range("E23").Value = 935.8685416337494
$ws.Range("F23").Value = 8.596431022649288
$ws.Range("G23").Value = 4.162330905306977
$ws.Range("H23").Value = 26.91712916011075
